$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row 2, column B: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the two "section header only" rows that had no data:
#  row 8 "grandes regiões" and row 5 "situação do domicílio"
# Delete row 8 first so row indices above it (row 5) remain valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
